$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 9 data rows (2-10) down by two rows (to 4-12), working
# from the bottom up so we never overwrite a row before it has been copied.
# Plain value copy (no Range.Insert) keeps the untouched data rows free of
# any inherited formatting.
for ($r = 10; $r -ge 2; $r--) {
    $dest = $r + 2
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Copy($ws.Cells.Item($dest, $c))
    }
}

# New row 2: Morning Run (2025-07-14)
$ws.Range("A2").Value = "Morning Run"
$ws.Range("B2").Value = 9.31
$ws.Range("C2").Value = 3983
$ws.Range("D2").Value = "Run"
$ws.Range("E2").Value = "2025-07-14T09:56:26Z"
$ws.Range("F2").Value = "07:08"
$ws.Range("G2").Value = 2.339
$ws.Range("H2").Value = 140.5

# New row 3: Afternoon Run (2025-07-12)
$ws.Range("A3").Value = "Afternoon Run"
$ws.Range("B3").Value = 18.41
$ws.Range("C3").Value = 5523
$ws.Range("D3").Value = "Run"
$ws.Range("E3").Value = "2025-07-12T16:48:18Z"
$ws.Range("F3").Value = "05:00"
$ws.Range("G3").Value = 3.334
$ws.Range("H3").Value = 168.3
